# DiffExcel V1.2: add Sheet4 (a copy-pattern of Sheet3's "same" block)
# after Sheet3, fill B2:F13 with "same", and make Sheet4 the active tab.

$wb = $excel.ActiveWorkbook

# Add a new worksheet and name it Sheet4.
$newSheet = $wb.Worksheets.Add()
$newSheet.Name = "Sheet4"

# Reposition it right after Sheet3 (Add() drops it before the active sheet).
$ws4 = $wb.Worksheets.Item("Sheet4")
$ws4.Move($null, $wb.Worksheets.Item("Sheet3"))

# Fill B2:F13 with the shared "same" text, matching Sheet3's pattern.
$ws4 = $wb.Worksheets.Item("Sheet4")
$ws4.Range("B2:F13").Value = "same"

# Make Sheet4 the active/selected tab, with F14 as the selected cell.
$ws4.Select()
$ws4.Range("F14").Select()
